# Applies the "sort banks by name" commit:
#  1. accounts sheet: swap sortOrder between "Citi Checking" (row2) and
#     "Citi Credit" (row4), and update their _history audit-log blobs
#     to reflect the resort events.
#  2. banks sheet: add a new bank "1st Source Bank" as row 4 (keeping
#     banks sorted alphabetically: Citibank, Chase, ... 1st Source Bank
#     was appended at the end of the banks table).
#  3. settings sheet: record the current sidebar width (176px) that was
#     captured as part of this UI change.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1) accounts sheet - reorder "Citi Checking" / "Citi Credit"
# ---------------------------------------------------------------------
$accounts = $wb.Worksheets.Item("accounts")

$accounts.Cells.Item(2, 3).NumberFormat = "@"
$accounts.Cells.Item(2, 3).Value = "fY7LCsIwEEV/RUaXFZqmSZNui6gI+gHiwmrU0Ec0TYVS+u8mVaEidDWXM3OGu2/BQOx78IC4hVkljJvlsRAQQyKNnCQ3ccpkeQUPTHPv8WqRbNbbpSUnlStt0TS9hBhhS8q6SIVDNj9lJdPcKkbXwgOtauMe9btMNO9QKW12+uycOeq6zusbIUIiFlIWEoKiT7vB5beqPxQ4jTgjlAcjAvoVmI8xo3hECP4EjkJGx4XDCw=="
$accounts.Cells.Item(2, 13).Value = 2

$accounts.Cells.Item(4, 3).NumberFormat = "@"
$accounts.Cells.Item(4, 3).Value = "jY5NC4JAEIb/SkwdDVzdL72FdugUSLfokLrFkrq1roGI/73VCoQoOs3M884D774DA6HrwA3CDha1MMOsjqWAECJp5CzSIpcGHDDtdYTJOt7solUSW5apQmkL5yeW+R6zpGrKVAyIEmzPu6xlWljP6EY4oFVjZHW2sc0uon0utdJmq/NBW6K+752xFCKEcUw5JgSxV8HJ57utNxUCygJOaOD9L3DX9zn1fwjoQwgQ5vS74Frh8AA="
$accounts.Cells.Item(4, 13).Value = 0

# ---------------------------------------------------------------------
# 2) banks sheet - append "1st Source Bank"
# ---------------------------------------------------------------------
$banks = $wb.Worksheets.Item("banks")

$banks.Cells.Item(4, 1).NumberFormat = "@"
$banks.Cells.Item(4, 1).Value = "a3918c067d67964dd150fcfec19e6b4da"
$banks.Cells.Item(4, 2).Value = 0
$banks.Cells.Item(4, 3).NumberFormat = "@"
$banks.Cells.Item(4, 3).Value = "bVTbkqJIEP0Vg9i3NmxEcMWJfkBQKAdoxUbs2tkHqKqGUi6OgFwm+t+30DZGdzcqMiozK+vkycPlr19cwU35PveTm/7i/shJ0e2pnxBuyg3zorfJyhMivZmfHrg+V5GA5auqGrCz/HKUpTFNyQBlCTv3MT6RPGc1Ai/1vEFvST4+yCnP0h/ppqQF6Ym8yNysLKLejKS43wN2TxyP+eGP1N0oDCLNCtIBMJeiLL161ybctDiVpM99UMzS/J9DoVuT7vwUsowGOvejZm5UFMd8+vzMIkwHmIY0zVn7juYzSk6dDbrKPlfm5PQ1L4uOfp5X2anD5z4/P/sXeYaSJI8nvMQzE76kunK7ScZhv/CnNPFD8pyfw6c6ib9VFBfRy3D8LSI0jIrOC/ycjMX+SpPFwKtCnGwbJMTnYM9TayNWYD95UqkSAlWZANdZrYcOdOduiEYwBYvZdn2ol+tGoeZmct62tmY28nw7X4fuEBpAm5cWyzuuSIG6zLDhVK90csYjPDJT1JqJ3MBGNlCyqHw93qNGftse0Jn1rU3ecR1VbnGCamsjHbC+pqsw63iEK8MWYDNjfOsStTz1DYdHWnY2GS5upJHVSGeUMJy9UllXDAr0KPY9nGGNzfWmVMsNCLEQM9xQBvs5q3MpMOAR7rAajEJ52SqhpSq1TTtDNw1CoNUH6MH2xuUr12IjzuFmxpPdLF5tljzcRbwp2C1qwO+7VyuREEVQl0fAyL/f5S+G9a2I1WHr64tDMEJjoM0qoioVMsLQerv4e0tTqlftX3cN/j6emCObJ14d3/GcmIITw9S6n6WBns0DHR4DvZIBtcRXLWQ6dJoobD9c9wft3lkeH6HhZKsNYNqAp5UqX3Due13mUO400u38fWe3TJvbbPReF5jEecDwgQCTy6Lgbp5ZEggSb3owCrw4J288RcKiRM2whbvl8X9rd/aRJC7j75ZW+98F7p+f4cTEWBfvnrT3dblhHAvfcw6B7j5y1GWWl6IgZbXJooVefQy6HokdB6nTvHvVQz1h87yqD9wkppl41/vh3bC2tvr72VTni4aa+P0Wd9/cRWPl5aX7E/z9Dw=="
$banks.Cells.Item(4, 4).NumberFormat = "@"
$banks.Cells.Item(4, 4).Value = "1st Source Bank"
$banks.Cells.Item(4, 5).NumberFormat = "@"
$banks.Cells.Item(4, 5).Value = "www.1stsourceonline.com"
$banks.Cells.Item(4, 6).NumberFormat = "@"
$banks.Cells.Item(4, 6).Value = "205 W. Jefferson`nSuite 404`nSouth Bend, IN 46601`nUSA"
$banks.Cells.Item(4, 7).NumberFormat = "@"
$banks.Cells.Item(4, 7).Value = ""
$banks.Cells.Item(4, 8).NumberFormat = "@"
$banks.Cells.Item(4, 8).Value = "data:image/svg+xml;width=16;height=16;base64,PD94bWwgdmVyc2lvbj0iMS4wIj8+CiAgICA8IURPQ1RZUEUgc3ZnIFBVQkxJQyAiLS8vVzNDLy9EVEQgU1ZHIDEuMS8vRU4iICJodHRwOi8vd3d3LnczLm9yZy9HcmFwaGljcy9TVkcvMS4xL0RURC9zdmcxMS5kdGQiPgogICAgPHN2ZyB4bWxucz0iaHR0cDovL3d3dy53My5vcmcvMjAwMC9zdmciIGhlaWdodD0iMTAwJSIgd2lkdGg9IjEwMCUiIHZpZXdCb3g9JzAgMCAxNiAxNic+CiAgICAgIDxkZWZzPgogICAgICAgIDxzdHlsZSB0eXBlPSJ0ZXh0L2NzcyI+CiAgICAgICAgICAuc2hhZG93IHsKICAgICAgICAgICAgdGV4dC1zaGFkb3c6IDBweCAwcHggMTBweCAjMDAwODsKICAgICAgICAgIH0KICAgICAgICA8L3N0eWxlPgogICAgICA8L2RlZnM+CiAgICAgIDxyZWN0IGZpbGw9IiM4ODgiIHg9IjAiIHk9IjAiIGhlaWdodD0iMTYiIHdpZHRoPSIxNiI+PC9yZWN0PgogICAgICA8dGV4dAogICAgICAgIGNsYXNzPSJzaGFkb3ciCiAgICAgICAgZmlsbD0iI2ZmZmZmZiIKICAgICAgICBmb250LWZhbWlseT0ic2Fucy1zZXJpZiIKICAgICAgICBmb250LXNpemU9IjUuMzMzMzMzMzMzMzMzMzMzIgogICAgICAgIHRleHQtYW5jaG9yPSJtaWRkbGUiCiAgICAgICAgZG9taW5hbnQtYmFzZWxpbmU9ImNlbnRyYWwiCiAgICAgICAgeD0iOCIKICAgICAgICB5PSI4IgogICAgICA+CiAgICAgICAgMVNCCiAgICAgIDwvdGV4dD4KICAgIDwvc3ZnPgogIA=="
$banks.Cells.Item(4, 9).Value = $true
$banks.Cells.Item(4, 10).NumberFormat = "@"
$banks.Cells.Item(4, 10).Value = "071212128"
$banks.Cells.Item(4, 11).NumberFormat = "@"
$banks.Cells.Item(4, 11).Value = "DI"
$banks.Cells.Item(4, 12).NumberFormat = "@"
$banks.Cells.Item(4, 12).Value = "https://ofxdi.diginsite.com/cmr/cmr.ofx"
$banks.Cells.Item(4, 13).NumberFormat = "@"
$banks.Cells.Item(4, 13).Value = ""
$banks.Cells.Item(4, 14).NumberFormat = "@"
$banks.Cells.Item(4, 14).Value = ""

# ---------------------------------------------------------------------
# 3) settings sheet - persist sidebar width setting
# ---------------------------------------------------------------------
$settings = $wb.Worksheets.Item("settings")

$settings.Cells.Item(2, 1).NumberFormat = "@"
$settings.Cells.Item(2, 1).Value = "sidebarWidth"
$settings.Cells.Item(2, 2).NumberFormat = "@"
$settings.Cells.Item(2, 2).Value = "176"

Write-Output "applied sort-banks-by-name edit"
